# Updates the crypto price ("D") and 1h-volume-change ("E") columns with
# refreshed values, row by row (row 1 is the header).
#
# Numeric-looking price strings (e.g. "1.00", "0.0820") are written with a
# leading apostrophe so Excel keeps them as text instead of normalising them
# into numbers (which would silently drop the trailing zero / change the
# stored type) - exactly what typing '1.00 into a cell does in the Excel UI.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '48.293.35'
$ws.Range('E2').Value = '  +1.91%  '
# Row 3
$ws.Range('D3').Value = '2.523.63'
$ws.Range('E3').Value = '  +0.79%  '
# Row 4
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.05%  '
# Row 5
$ws.Range('D5').Value = '''323.20'
# Row 6
$ws.Range('D6').Value = '''109.19'
$ws.Range('E6').Value = '  -0.73%  '
# Row 7
$ws.Range('E7').Value = '  +0.47%  '
# Row 8
$ws.Range('E8').Value = '  -0.03%  '
# Row 9
$ws.Range('E9').Value = '  +4.21%  '
# Row 10
$ws.Range('D10').Value = '''40.36'
$ws.Range('E10').Value = '  +1.99%  '
# Row 11
$ws.Range('D11').Value = '''20.04'
$ws.Range('E11').Value = '  +8.09%  '
# Row 12
$ws.Range('D12').Value = '''0.0820'
$ws.Range('E12').Value = '  +0.49%  '
# Row 13
$ws.Range('E13').Value = '  +0.79%  '
# Row 14
$ws.Range('E14').Value = '  +0.43%  '
# Row 15
$ws.Range('D15').Value = '2.915.85'
$ws.Range('E15').Value = '  +0.68%  '
# Row 16
$ws.Range('D16').Value = '2.522.05'
$ws.Range('E16').Value = '  +0.46%  '
# Row 17
$ws.Range('E17').Value = '  -0.24%  '
# Row 18
$ws.Range('D18').Value = '48.172.96'
$ws.Range('E18').Value = '  +1.77%  '
# Row 19
$ws.Range('D19').Value = '''13.27'
$ws.Range('E19').Value = '  +3.08%  '
# Row 20
$ws.Range('E20').Value = '  -0.44%  '
# Row 21
$ws.Range('D21').Value = '0.0₃0945'
$ws.Range('E21').Value = '  +0.26%  '
# Row 22
$ws.Range('D22').Value = '''2.72'
$ws.Range('E22').Value = '  +0.55%  '
# Row 23
$ws.Range('D23').Value = '''72.54'
$ws.Range('E23').Value = '  +2.80%  '
# Row 24
$ws.Range('D24').Value = '''268.29'
# Row 25
$ws.Range('E25').Value = '  -0.94%  '
# Row 26
$ws.Range('E26').Value = '  +0.26%  '
# Row 27
$ws.Range('E27').Value = '  +0.07%  '
# Row 28
$ws.Range('E28').Value = '  +0.56%  '
# Row 29
$ws.Range('D29').Value = '''10.21'
$ws.Range('E29').Value = '  +1.32%  '
# Row 30
$ws.Range('E30').Value = '  +4.96%  '
# Row 31
$ws.Range('D31').Value = '''35.13'
$ws.Range('E31').Value = '  -0.90%  '
# Row 32
$ws.Range('E32').Value = '  -0.22%  '
# Row 33
$ws.Range('D33').Value = '''19.99'
$ws.Range('E33').Value = '  +0.00%  '
# Row 34
$ws.Range('E34').Value = '  -0.50%  '
# Row 35
$ws.Range('E35').Value = '  -0.04%  '
# Row 36
$ws.Range('E36').Value = '  -0.65%  '
# Row 37
$ws.Range('E37').Value = '  -0.67%  '
# Row 38
$ws.Range('E38').Value = '  +0.47%  '
# Row 39
$ws.Range('D39').Value = '''3.00'
$ws.Range('E39').Value = '  +0.09%  '
# Row 40
$ws.Range('E40').Value = '  +0.16%  '
# Row 41
$ws.Range('D41').Value = '''22.31'
$ws.Range('E41').Value = '  +4.94%  '
# Row 42
$ws.Range('E42').Value = '  -1.41%  '
# Row 43
$ws.Range('D43').Value = '''118.38'
$ws.Range('E43').Value = '  -2.67%  '
# Row 44
$ws.Range('E44').Value = '  -0.03%  '
# Row 45
$ws.Range('D45').Value = '2.001.93'
$ws.Range('E45').Value = '  +0.03%  '
# Row 46
$ws.Range('E46').Value = '  +0.03%  '
# Row 47
$ws.Range('E47').Value = '  +6.49%  '
# Row 48
$ws.Range('E48').Value = '  -2.02%  '
# Row 49
$ws.Range('D49').Value = '''9.09'
$ws.Range('E49').Value = '  +0.39%  '
# Row 51
$ws.Range('D51').Value = '''80.64'
$ws.Range('E51').Value = '  +3.01%  '
